# Assets.xlsx maintenance pass:
#  - stamp every vendor sheet with a "Vendor ID" / "MPulse ID" header row
#    (feeds the upcoming PartsTable class / archive lookups)
#  - leave the whole header row selected on each tab, the way it was left
#    after typing the headers in and hitting Ctrl+Space
#  - auto-size the id columns on the first sheet
#  - finish up on the EverPower tab, which becomes the active sheet

$wb = $excel.ActiveWorkbook

$vendorIdHeader = "Vendor ID"
$mpulseIdHeader = "MPulse ID"

$sheetCount = $wb.Worksheets.Count
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $ws.Range("A1").Value = $vendorIdHeader
    $ws.Range("B1").Value = $mpulseIdHeader

    # Select the full header row (as if the row header had been clicked).
    $ws.Rows.Item(1).Select() | Out-Null
}

# The first sheet also got its new columns auto-fit to the header text.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns("A:B").AutoFit() | Out-Null
$ws1.Rows.Item(1).Select() | Out-Null

# EverPower (last sheet) is the one left active/visible.
$lastSheet = $wb.Worksheets.Item($sheetCount)
$lastSheet.Activate()
$lastSheet.Rows.Item(1).Select() | Out-Null
